$wb = $excel.ActiveWorkbook

# --- Update the scaling values (order matters: it determines the order new
# shared strings are appended in the saved workbook) ---
$wb.Worksheets.Item("VcpuValueForHScaling").Range("B2").Value = "4"
$wb.Worksheets.Item("RAMValueForHScaling").Range("B2").Value = "4"
$wb.Worksheets.Item("DiskSizeForHScaling").Range("B2").Value = "36"
$wb.Worksheets.Item("VLANName").Range("B2").Value = "10.150.43.0 - Test_01"
$wb.Worksheets.Item("VMNamesForH").Range("B2").Value = "test2"

# --- UserInfo: update the role/user entries that used to hold lead2/lead3 ---
$uiSheet = $wb.Worksheets.Item("UserInfo")
$uiSheet.Range("C3").Value = "test1"
$uiSheet.Range("C4").Value = "test2"

# --- Storage page values (StoragePathName sheet, and the copy shown on OSName) ---
$wb.Worksheets.Item("StoragePathName").Range("B2").Value = "Storage3_2025 - 122GB"
$wb.Worksheets.Item("OSName").Range("C2").Value = "Storage3_2025 - 122GB"

# --- View/selection state ---
# Leave a lingering selection on the UserInfo sheet at C4 (not the active tab).
$uiSheet.Activate()
$uiSheet.Range("C4").Select()

# Finish with the OSName sheet active/selected - this is the sheet the author
# was working on ("working on Storage page").
$osSheet = $wb.Worksheets.Item("OSName")
$osSheet.Activate()
$osSheet.Range("C2").Select()
